# Added LRU results for mcf (and completed the milc block that the LRU
# simulation run also produced) on both the Config1 and Config2 sheets.
#
# Rows 55 (mcf / LRU) and 59-62 (milc / LRU, SRRIP, Hawkeye, OPTGen) had
# placeholder #DIV/0! errors in H (IPC) and I (MPKI) because the raw
# simulation counters in C:G were never filled in. This fills in the
# counters that the simulator produced, which lets the existing shared
# formulas in H/I (and the local formulas in G62/J62) recalculate to real
# numbers instead of #DIV/0!.

$wb = $excel.ActiveWorkbook

function Fill-Row {
    param(
        $ws,
        [int]$row,
        [double]$c,
        [double]$d,
        [double]$e,
        [double]$f,
        [double]$g = [double]::NaN,
        [bool]$isOptGen = $false
    )

    $ws.Cells.Item($row, 3).Value = $c   # C - Total Instructions
    $ws.Cells.Item($row, 4).Value = $d   # D - Total Cycles
    $ws.Cells.Item($row, 5).Value = $e   # E - Total Access
    $ws.Cells.Item($row, 6).Value = $f   # F - Total Hit

    if ($isOptGen) {
        # OPTGen rows compute Total Miss (G) and Hit rate (J) from E/F,
        # matching the pattern already used elsewhere (e.g. rows 50 & 54).
        $ws.Cells.Item($row, 7).Formula = "=E$row-F$row"
        $ws.Cells.Item($row, 10).Formula = "=F$row/E$row"
    } else {
        $ws.Cells.Item($row, 7).Value = $g   # G - Total Miss
    }

    # H - IPC, I - MPKI: same shared formulas used throughout the table.
    $ws.Cells.Item($row, 8).Formula = "=(C$row/D$row)"
    $ws.Cells.Item($row, 9).Formula = "=G$row/(C$row/1000)"
}

# ---- Config1 sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config1")

Fill-Row $ws1 55 50000001 505789892 5694538 2113376 3581162
Fill-Row $ws1 59 50000002 121522048 1510289 325673  1184616
Fill-Row $ws1 60 50000002 122348912 1510289 299898  1210391
Fill-Row $ws1 61 50000002 120220593 1510289 38749   1471540
Fill-Row $ws1 62 50000002 120220593 32737   571     0 $true

# ---- Config2 sheet -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config2")

Fill-Row $ws2 55 50000000 437451603 9525395 4305218 5220177
Fill-Row $ws2 59 50000002 93003755  1856592 325679  1530913
Fill-Row $ws2 60 50000002 93552246  1856586 300007  1556579
Fill-Row $ws2 61 50000002 92679825  1856587 32820   1823767
Fill-Row $ws2 62 50000002 92679825  17406   170     0 $true

# ---- View-state: Config1 becomes the active sheet, both sheets' current
# selection moves from C59 to C56, and the view scrolls back up a bit.
$ws2.Range("C56").Select()
$excel.ActiveWindow.ScrollRow = 48

$ws1.Activate()
$ws1.Range("C56").Select()
$excel.ActiveWindow.ScrollRow = 46
